# Replace the multiplication-fact answers in the table cells with new
# regenerated problems/answers, per the commit "Update master to output
# generated at 4250d90".
$d = $word.ActiveDocument

$d.Content.Find.Execute("977×7=6839", $true, $true, $false, $false, $false, $true, 1, $false, "511×7=3577", 2) | Out-Null
$d.Content.Find.Execute("409×3=1227", $true, $true, $false, $false, $false, $true, 1, $false, "135×8=1080", 2) | Out-Null
$d.Content.Find.Execute("218×4=872", $true, $true, $false, $false, $false, $true, 1, $false, "931×5=4655", 2) | Out-Null
$d.Content.Find.Execute("723×5=3615", $true, $true, $false, $false, $false, $true, 1, $false, "206×3=618", 2) | Out-Null
$d.Content.Find.Execute("629×6=3774", $true, $true, $false, $false, $false, $true, 1, $false, "110×4=440", 2) | Out-Null
$d.Content.Find.Execute("969×7=6783", $true, $true, $false, $false, $false, $true, 1, $false, "754×6=4524", 2) | Out-Null
$d.Content.Find.Execute("675×4=2700", $true, $true, $false, $false, $false, $true, 1, $false, "849×3=2547", 2) | Out-Null
$d.Content.Find.Execute("195×9=1755", $true, $true, $false, $false, $false, $true, 1, $false, "341×7=2387", 2) | Out-Null
$d.Content.Find.Execute("652×9=5868", $true, $true, $false, $false, $false, $true, 1, $false, "370×6=2220", 2) | Out-Null
$d.Content.Find.Execute("178×5=890", $true, $true, $false, $false, $false, $true, 1, $false, "951×8=7608", 2) | Out-Null
$d.Content.Find.Execute("948×8=7584", $true, $true, $false, $false, $false, $true, 1, $false, "752×5=3760", 2) | Out-Null
$d.Content.Find.Execute("881×5=4405", $true, $true, $false, $false, $false, $true, 1, $false, "246×6=1476", 2) | Out-Null
$d.Content.Find.Execute("701×2=1402", $true, $true, $false, $false, $false, $true, 1, $false, "404×8=3232", 2) | Out-Null
$d.Content.Find.Execute("626×6=3756", $true, $true, $false, $false, $false, $true, 1, $false, "571×9=5139", 2) | Out-Null
$d.Content.Find.Execute("366×9=3294", $true, $true, $false, $false, $false, $true, 1, $false, "417×2=834", 2) | Out-Null
$d.Content.Find.Execute("862×2=1724", $true, $true, $false, $false, $false, $true, 1, $false, "920×4=3680", 2) | Out-Null
$d.Content.Find.Execute("919×6=5514", $true, $true, $false, $false, $false, $true, 1, $false, "885×2=1770", 2) | Out-Null
$d.Content.Find.Execute("355×8=2840", $true, $true, $false, $false, $false, $true, 1, $false, "166×4=664", 2) | Out-Null
$d.Content.Find.Execute("851×7=5957", $true, $true, $false, $false, $false, $true, 1, $false, "453×5=2265", 2) | Out-Null
$d.Content.Find.Execute("152×7=1064", $true, $true, $false, $false, $false, $true, 1, $false, "304×2=608", 2) | Out-Null
$d.Content.Find.Execute("921×5=4605", $true, $true, $false, $false, $false, $true, 1, $false, "173×9=1557", 2) | Out-Null
$d.Content.Find.Execute("372×7=2604", $true, $true, $false, $false, $false, $true, 1, $false, "221×7=1547", 2) | Out-Null
$d.Content.Find.Execute("494×3=1482", $true, $true, $false, $false, $false, $true, 1, $false, "900×6=5400", 2) | Out-Null
$d.Content.Find.Execute("471×2=942", $true, $true, $false, $false, $false, $true, 1, $false, "765×3=2295", 2) | Out-Null
$d.Content.Find.Execute("857×9=7713", $true, $true, $false, $false, $false, $true, 1, $false, "194×3=582", 2) | Out-Null
